$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.497.82"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.154.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.12%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.70%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.149.63"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.22%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.09%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.58%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.18%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.27%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.52"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.80%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.672.45"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.14%  "

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.93%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.447.15"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.155.93"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.84"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "476.19"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.64"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.722"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.32%  "

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.34%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.77"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.86"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +2.42%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.05%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.80"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -3.17%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.20%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +8.53%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.77%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.58%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.43"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.65"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.70%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.82%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.71"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0748"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.10"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "451.94"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0395"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.77%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.84%  "

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.846.46"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.267"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.98%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.20%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.45"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.45"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.29%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.53%  "

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.03%  "
